{"js": "// Update the header date and all table answer cells to the new values.\n\n// 1) Header date paragraph: \"2025-12-09 Tuesday\" -> \"2025-12-10 Wednesday\"\nconst dateResults = context.document.body.search(\"2025-12-09 Tuesday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2025-12-10 Wednesday\", Word.InsertLocation.replace);\n}\n\n// 2) Table cells: replace the division problems row by row, in document order.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = [\n  [\"813\u00f72=406, 1\", \"502\u00f77=71, 5\", \"423\u00f73=141, 0\", \"445\u00f75=89, 0\", \"490\u00f76=81, 4\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"364\u00f74=91, 0\", \"905\u00f76=150, 5\", \"305\u00f79=33, 8\", \"184\u00f76=30, 4\", \"812\u00f77=116, 0\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"472\u00f74=118, 0\", \"403\u00f75=80, 3\", \"457\u00f78=57, 1\", \"378\u00f79=42, 0\", \"966\u00f74=241, 2\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"970\u00f77=138, 4\", \"479\u00f77=68, 3\", \"823\u00f74=205, 3\", \"176\u00f73=58, 2\", \"452\u00f77=64, 4\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"126\u00f77=18, 0\", \"160\u00f75=32, 0\", \"335\u00f72=167, 1\", \"222\u00f72=111, 0\", \"221\u00f75=44, 1\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n];\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Header date: \"2025-12-09 Tuesday\" -> \"2025-12-10 Wednesday\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"2025-12-09 Tuesday\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"2025-12-10 Wednesday\"\n$find.Execute([ref]\"2025-12-09 Tuesday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2025-12-10 Wednesday\", 2)\n\n# 2) Table cells: replace the division problems, row by row (only every 4th\n#    row of the 20-row table actually holds text; the rest are blank spacer\n#    rows), left to right within each row.\n$tbl = $d.Tables.Item(1)\n$values = @(\n  @(\"813\u00f72=406, 1\", \"502\u00f77=71, 5\", \"423\u00f73=141, 0\", \"445\u00f75=89, 0\", \"490\u00f76=81, 4\"),\n  @(\"364\u00f74=91, 0\", \"905\u00f76=150, 5\", \"305\u00f79=33, 8\", \"184\u00f76=30, 4\", \"812\u00f77=116, 0\"),\n  @(\"472\u00f74=118, 0\", \"403\u00f75=80, 3\", \"457\u00f78=57, 1\", \"378\u00f79=42, 0\", \"966\u00f74=241, 2\"),\n  @(\"970\u00f77=138, 4\", \"479\u00f77=68, 3\", \"823\u00f74=205, 3\", \"176\u00f73=58, 2\", \"452\u00f77=64, 4\"),\n  @(\"126\u00f77=18, 0\", \"160\u00f75=32, 0\", \"335\u00f72=167, 1\", \"222\u00f72=111, 0\", \"221\u00f75=44, 1\")\n)\n\n$rowIndex = 1\nforeach ($rowVals in $values) {\n  $colIndex = 1\n  foreach ($val in $rowVals) {\n    $cell = $tbl.Cell($rowIndex, $colIndex)\n    $cell.Range.Text = $val\n    $colIndex = $colIndex + 1\n  }\n  $rowIndex = $rowIndex + 4\n}\n"}
